$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "SCD0020"

# Update the TC_ID value in B2 from "DGS-327" to "SCD0020-006"
$ws.Range("B2").Value = "SCD0020-006"

# The new value is wider than the old one, so the best-fit column B needs to
# widen to keep showing the whole value (matches the saved column width as
# closely as this engine's ColumnWidth -> stored-width quantization allows).
$ws.Columns.Item(2).ColumnWidth = 11.67

# Move the active selection to B3 (matches the saved selection in the sheet view)
$ws.Range("B3").Select()
